# Update loading_percent values on Sheet1 for the "380 kV" case.
# Each data row (rows 2-25) gets new values in columns B, C, E, F, G, H, J, N.
# Columns A, D, I, K, L, M, O are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row, B, C, E, F, G, H, J, N
$data = @(
    @(2,14.83835186183086,9.23881429251133,22.4311056418522,41.31237264083453,30.45625978056599,14.18117724149123,7.907545675285167,16.41873126027958),
    @(3,14.21895140864909,8.653321671085424,22.24519652872732,41.01159636432503,30.16374241973622,14.2131078255066,7.934231685557175,16.48619248077645),
    @(4,13.82756096988738,8.271744426313541,22.13486707949318,40.83918225546569,29.99843427120287,14.2367525506876,7.951934677526477,16.52964571338262),
    @(5,13.6655312040554,8.110695898619259,22.09090639268269,40.77206648729512,29.93473730037935,14.24739846758438,7.959479813818358,16.54786571664449),
    @(6,13.63848108017719,8.083618624786377,22.0836682221131,40.76111352118131,29.92438386734174,14.24922708763863,7.960752661846166,16.55092213143395),
    @(7,13.82538569122558,8.269594949300828,22.13427011340404,40.83826430221803,29.99756029515122,14.23689204094932,7.952035093990653,16.52988935776443),
    @(8,14.62722234460613,9.041518495088923,22.36623540414017,41.20615791402141,30.35247842304625,14.19134553619074,7.916473336963582,16.44157112906894),
    @(9,16.10221970692436,10.38023309081922,22.84962840895919,42.02194287715781,31.15826510220174,14.1342970729259,7.85720641652514,16.28442965013956),
    @(10,17.11567299725478,11.25755666207094,23.21965259701924,42.67446116367862,31.81167561817378,14.11233347833209,7.820065726546447,16.17866325740617),
    @(11,17.55966509660999,11.63374806786294,23.39066035967643,42.98182429779027,32.12093221708577,14.10672551998571,7.80456429881087,16.13262964729274),
    @(12,17.72522174518693,11.77292030726875,23.45575442565703,43.09963677736697,32.23964875946537,14.10523583401518,7.798895188101946,16.11549549844681),
    @(13,17.68968240016191,11.74309291771262,23.4417210352295,43.07420208192359,32.21401147374933,14.10552841780138,7.800107187781172,16.11917242096508),
    @(14,17.573337765173,11.64526358821735,23.39600920990109,42.99148887299254,32.13066757767525,14.10659023873099,7.804093867676015,16.13121405011017),
    @(15,17.50173479653967,11.58491305659016,23.36805185694283,42.94100684401783,32.07982278435949,14.10732329002156,7.806562007384244,16.13862863646704),
    @(16,17.08630327852962,11.23251160981999,23.20852666507696,42.65457822453176,31.79169659245754,14.11278849422235,7.821106874133668,16.18171340383635),
    @(17,16.82699257306655,11.01046491211897,23.11131563009564,42.48149491085589,31.61793404053922,14.11726660940963,7.830387069154686,16.20867631009784),
    @(18,16.67624685449589,10.88059264142109,23.05565778265221,42.38293884011375,31.51912898973356,14.12025470006814,7.835855978077387,16.22438052277007),
    @(19,16.62493674122281,10.83624936427266,23.03685832382657,42.34974345603318,31.48587444723439,14.12133712903534,7.837730173535159,16.22973137976226),
    @(20,16.85476288880599,11.03432531630635,23.12163783142391,42.49981741954884,31.63631438205074,14.11674719756794,7.829385595856975,16.20578580151939),
    @(21,17.60758174825295,11.67408745001461,23.40942711941662,43.01574594776076,32.15510504162098,14.10626112507476,7.802917425737925,16.12766906051566),
    @(22,18.08454684126863,12.07307718719632,23.59945648895455,43.36117095840618,32.50347805993211,14.10310423710951,7.786790427732253,16.0783503488261),
    @(23,17.83139464678122,11.86187559995633,23.49787259950886,43.17608911001135,32.31673324493497,14.10444983354669,7.795290358811837,16.10451434394107),
    @(24,16.842213104098,11.02354492665703,23.11697044388634,42.49153084141113,31.62800121486159,14.11698073570161,7.829837945821098,16.2070919683768),
    @(25,15.71483962791932,10.03677268932705,22.71605935419419,41.79158037466685,30.92906491082926,14.14624546644998,7.872116798070413,16.32523264807015)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Range("B$row").Value = $entry[1]
    $ws.Range("C$row").Value = $entry[2]
    $ws.Range("E$row").Value = $entry[3]
    $ws.Range("F$row").Value = $entry[4]
    $ws.Range("G$row").Value = $entry[5]
    $ws.Range("H$row").Value = $entry[6]
    $ws.Range("J$row").Value = $entry[7]
    $ws.Range("N$row").Value = $entry[8]
}
